$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column header B1
$ws.Range("B1").Value = "n_contornos"

# Add new headers E1 and F1
$ws.Range("E1").Value = "min_q"
$ws.Range("F1").Value = "max_q"

# Match the header formatting used by the existing header row (D1)
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data for columns E (min_q) and F (max_q), rows 2-18
$minVals = @{
    2  = 0.00002304
    3  = 0.00002056
    4  = 0.00000776
    5  = 0
    6  = 0.00003456
    7  = 0.0000324
    8  = 0.00005072
    9  = 0.00005552
    10 = 0.00003624
    11 = 0.00003304
    12 = 0.00005008
    13 = 0.00005424
    14 = 0.00001848
    15 = 0.00002344
    16 = 0.00001616
    17 = 0.00001736
    18 = 0.00001496
}

$maxVals = @{
    2  = 0.01538384
    3  = 0.00589112
    4  = 0.0084864
    5  = 0.01146816
    6  = 0.01113464
    7  = 0.0109924
    8  = 0.01420992
    9  = 0.0091236
    10 = 0.00834336
    11 = 0.00639712
    12 = 0.00914432
    13 = 0.00627856
    14 = 0.00869776
    15 = 0.00897856
    16 = 0.01161448
    17 = 0.00934504
    18 = 0.00642688
}

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 5).Value = $minVals[$r]
    $ws.Cells.Item($r, 6).Value = $maxVals[$r]
}
